$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet's tab / sheet name
$ws.Name = "UniformF"

# Add new row 16 data
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16:M16").Value = 1

# Copy the formatting (style) of row 15 into row 16
$ws.Range("A15:M15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
